# "Common: Improved inventory stuff"
#
# Inserts a new translation row for "common.inventory.code.tooltip" at the
# top of the (already alphabetically-ish curated) "Translations - Common"
# sheet list, right after the existing "common.job.*" style entries, pushing
# every following row down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Translations - Common")

# Insert a new row at row 11 - every row from 11 downward shifts to 12+,
# carrying along its original formatting (row height, styles, ...).
$ws.Rows.Item(11).Insert()

# Populate the freshly inserted row with the new translation entry.
$ws.Range("A11").Value = "cs"
$ws.Range("B11").Value = "common.inventory.code.tooltip"
$ws.Range("C11").Value = "Kód v inventáři pro rozlišení produktů stejného typu."

# Re-apply the sheet's sort (it covers the whole translation table below the
# header) so the sort cache/range grows to match the now-larger data range.
$sort = $ws.Sort
$sort.SortFields.Clear()
$key = $ws.Range("B2:B57")
$sort.SortFields.Add($key) | Out-Null
$sortRange = $ws.Range("A2:C57")
$sort.SetRange($sortRange)
$sort.Header = 0
$sort.Apply()

# Restore focus/selection like the author left it.
$ws.Activate()
$ws.Range("B51").Select()
